$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.064198653302102
$ws.Cells.Item(2, 4).Value = 1.061752694166239
$ws.Cells.Item(2, 5).Value = 1.069076882600384
$ws.Cells.Item(2, 6).Value = 1.078871629810596
$ws.Cells.Item(2, 9).Value = 1.055988313694664
$ws.Cells.Item(2, 10).Value = 1.069160123005274
$ws.Cells.Item(2, 11).Value = 1.06447590068766
$ws.Cells.Item(2, 12).Value = 1.071780355622972
$ws.Cells.Item(2, 13).Value = 1.081549165381923

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.065430994089893
$ws.Cells.Item(3, 4).Value = 1.062697770233633
$ws.Cells.Item(3, 5).Value = 1.070189976582276
$ws.Cells.Item(3, 6).Value = 1.080090757819482
$ws.Cells.Item(3, 9).Value = 1.056427830243159
$ws.Cells.Item(3, 10).Value = 1.070046242116564
$ws.Cells.Item(3, 11).Value = 1.065235330747638
$ws.Cells.Item(3, 12).Value = 1.072708813494197
$ws.Cells.Item(3, 13).Value = 1.082585280450773

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.066227995820214
$ws.Cells.Item(4, 4).Value = 1.063308840469845
$ws.Cells.Item(4, 5).Value = 1.07091003792658
$ws.Cells.Item(4, 6).Value = 1.080879649755083
$ws.Cells.Item(4, 9).Value = 1.056710730055846
$ws.Cells.Item(4, 10).Value = 1.070618669374012
$ws.Cells.Item(4, 11).Value = 1.065725650571289
$ws.Cells.Item(4, 12).Value = 1.073308810367187
$ws.Cells.Item(4, 13).Value = 1.083255176655956

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.066562960649962
$ws.Cells.Item(5, 4).Value = 1.063565626036134
$ws.Cells.Item(5, 5).Value = 1.071212708961712
$ws.Cells.Item(5, 6).Value = 1.081211309921249
$ws.Cells.Item(5, 9).Value = 1.056829303615877
$ws.Cells.Item(5, 10).Value = 1.070859091538858
$ws.Cells.Item(5, 11).Value = 1.065931523095783
$ws.Cells.Item(5, 12).Value = 1.073560864351585
$ws.Cells.Item(5, 13).Value = 1.083536673403854

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.066619197232511
$ws.Cells.Item(6, 4).Value = 1.063608735139855
$ws.Cells.Item(6, 5).Value = 1.071263526297082
$ws.Cells.Item(6, 6).Value = 1.081266997758559
$ws.Cells.Item(6, 9).Value = 1.056849191696733
$ws.Cells.Item(6, 10).Value = 1.070899446252187
$ws.Cells.Item(6, 11).Value = 1.065966074906923
$ws.Cells.Item(6, 12).Value = 1.073603174535058
$ws.Cells.Item(6, 13).Value = 1.083583930507927

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.066232472005882
$ws.Cells.Item(7, 4).Value = 1.063312272075578
$ws.Cells.Item(7, 5).Value = 1.070914082397871
$ws.Cells.Item(7, 6).Value = 1.08088408137303
$ws.Cells.Item(7, 9).Value = 1.05671231584536
$ws.Cells.Item(7, 10).Value = 1.070621882793982
$ws.Cells.Item(7, 11).Value = 1.065728402462085
$ws.Cells.Item(7, 12).Value = 1.073312179048834
$ws.Cells.Item(7, 13).Value = 1.083258938528792

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.064615213917317
$ws.Cells.Item(8, 4).Value = 1.062072182364863
$ws.Cells.Item(8, 5).Value = 1.069453097115555
$ws.Cells.Item(8, 6).Value = 1.079283633699152
$ws.Cells.Item(8, 9).Value = 1.056137160777981
$ws.Cells.Item(8, 10).Value = 1.069459788893503
$ws.Cells.Item(8, 11).Value = 1.064732778404908
$ws.Cells.Item(8, 12).Value = 1.072094294141568
$ws.Cells.Item(8, 13).Value = 1.081899437965073

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.06176219230467
$ws.Cells.Item(9, 4).Value = 1.059883435266897
$ws.Cells.Item(9, 5).Value = 1.066877170374092
$ws.Cells.Item(9, 6).Value = 1.076463628233203
$ws.Cells.Item(9, 9).Value = 1.055112163930701
$ws.Cells.Item(9, 10).Value = 1.067404682556591
$ws.Cells.Item(9, 11).Value = 1.062970021233544
$ws.Cells.Item(9, 12).Value = 1.069942206946031
$ws.Cells.Item(9, 13).Value = 1.079499622537262

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.059857875839255
$ws.Cells.Item(10, 4).Value = 1.058421812081613
$ws.Cells.Item(10, 5).Value = 1.06515879068591
$ws.Cells.Item(10, 6).Value = 1.074583640562616
$ws.Cells.Item(10, 9).Value = 1.054421048199815
$ws.Cells.Item(10, 10).Value = 1.066029578861576
$ws.Cells.Item(10, 11).Value = 1.061789167546999
$ws.Cells.Item(10, 12).Value = 1.068503344111922
$ws.Cells.Item(10, 13).Value = 1.077896820993713

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.059032706522831
$ws.Cells.Item(11, 4).Value = 1.057788313359612
$ws.Cells.Item(11, 5).Value = 1.064414431103091
$ws.Cells.Item(11, 6).Value = 1.073769563612397
$ws.Cells.Item(11, 9).Value = 1.054119928966758
$ws.Cells.Item(11, 10).Value = 1.065432928899487
$ws.Cells.Item(11, 11).Value = 1.061276479869945
$ws.Cells.Item(11, 12).Value = 1.067879298075045
$ws.Cells.Item(11, 13).Value = 1.077202073679487

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.058726110640815
$ws.Cells.Item(12, 4).Value = 1.057552911116605
$ws.Cells.Item(12, 5).Value = 1.064137897236373
$ws.Cells.Item(12, 6).Value = 1.07346717241756
$ws.Cells.Item(12, 9).Value = 1.054007799083726
$ws.Cells.Item(12, 10).Value = 1.065211121181614
$ws.Cells.Item(12, 11).Value = 1.061085837401708
$ws.Cells.Item(12, 12).Value = 1.067647346103436
$ws.Cells.Item(12, 13).Value = 1.076943902803432

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.058791880675441
$ws.Cells.Item(13, 4).Value = 1.057603409912322
$ws.Cells.Item(13, 5).Value = 1.064197216784531
$ws.Cells.Item(13, 6).Value = 1.073532036698313
$ws.Cells.Item(13, 9).Value = 1.054031864048638
$ws.Cells.Item(13, 10).Value = 1.065258708108628
$ws.Cells.Item(13, 11).Value = 1.06112674027062
$ws.Cells.Item(13, 12).Value = 1.067697107565047
$ws.Cells.Item(13, 13).Value = 1.076999286379454

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.059007365090337
$ws.Cells.Item(14, 4).Value = 1.057768856834886
$ws.Cells.Item(14, 5).Value = 1.064391573668882
$ws.Cells.Item(14, 6).Value = 1.073744568013185
$ws.Cells.Item(14, 9).Value = 1.054110666010181
$ws.Cells.Item(14, 10).Value = 1.065414597996645
$ws.Cells.Item(14, 11).Value = 1.061260725544837
$ws.Cells.Item(14, 12).Value = 1.067860127994864
$ws.Cells.Item(14, 13).Value = 1.077180735466388

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.05914011996915
$ws.Cells.Item(15, 4).Value = 1.057870781826546
$ws.Cells.Item(15, 5).Value = 1.064511317259958
$ws.Cells.Item(15, 6).Value = 1.073875514605491
$ws.Cells.Item(15, 9).Value = 1.054159181258383
$ws.Cells.Item(15, 10).Value = 1.065510622282771
$ws.Cells.Item(15, 11).Value = 1.061343250759683
$ws.Cells.Item(15, 12).Value = 1.067960549866282
$ws.Cells.Item(15, 13).Value = 1.07729251749477

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.059912626974161
$ws.Cells.Item(16, 4).Value = 1.058463842402173
$ws.Cells.Item(16, 5).Value = 1.065208185100226
$ws.Cells.Item(16, 6).Value = 1.074637667344712
$ws.Cells.Item(16, 9).Value = 1.054440993175613
$ws.Cells.Item(16, 10).Value = 1.066069150666266
$ws.Cells.Item(16, 11).Value = 1.061823163935975
$ws.Cells.Item(16, 12).Value = 1.0685447385724
$ws.Cells.Item(16, 13).Value = 1.0779429136617

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.060397040772176
$ws.Cells.Item(17, 4).Value = 1.058835690246593
$ws.Cells.Item(17, 5).Value = 1.065645232855506
$ws.Cells.Item(17, 6).Value = 1.075115736014019
$ws.Cells.Item(17, 9).Value = 1.054617267167105
$ws.Cells.Item(17, 10).Value = 1.066419172520121
$ws.Cells.Item(17, 11).Value = 1.062123832829888
$ws.Cells.Item(17, 12).Value = 1.068910913248003
$ws.Cells.Item(17, 13).Value = 1.078350694789242

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.060679534639767
$ws.Cells.Item(18, 4).Value = 1.059052524414765
$ws.Cells.Item(18, 5).Value = 1.065900127590292
$ws.Cells.Item(18, 6).Value = 1.075394582652771
$ws.Cells.Item(18, 9).Value = 1.054719905266353
$ws.Cells.Item(18, 10).Value = 1.066623216668643
$ws.Cells.Item(18, 11).Value = 1.062299075808164
$ws.Cells.Item(18, 12).Value = 1.06912439952358
$ws.Cells.Item(18, 13).Value = 1.078588476939238

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.060775848253016
$ws.Cells.Item(19, 4).Value = 1.059126449404728
$ws.Cells.Item(19, 5).Value = 1.065987035448774
$ws.Cells.Item(19, 6).Value = 1.075489661756233
$ws.Cells.Item(19, 9).Value = 1.054754871786513
$ws.Cells.Item(19, 10).Value = 1.066692770547775
$ws.Cells.Item(19, 11).Value = 1.062358806775725
$ws.Cells.Item(19, 12).Value = 1.069197176394004
$ws.Cells.Item(19, 13).Value = 1.078669542750175

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.060345073623347
$ws.Cells.Item(20, 4).Value = 1.058795800543388
$ws.Cells.Item(20, 5).Value = 1.065598344661535
$ws.Cells.Item(20, 6).Value = 1.075064444069987
$ws.Cells.Item(20, 9).Value = 1.054598373215115
$ws.Cells.Item(20, 10).Value = 1.066381630674464
$ws.Cells.Item(20, 11).Value = 1.062091587579083
$ws.Cells.Item(20, 12).Value = 1.068871636241178
$ws.Cells.Item(20, 13).Value = 1.078306950957051

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.058943912843902
$ws.Cells.Item(21, 4).Value = 1.057720139408693
$ws.Cells.Item(21, 5).Value = 1.064334341685086
$ws.Cells.Item(21, 6).Value = 1.073681983043346
$ws.Cells.Item(21, 9).Value = 1.054087468556801
$ws.Cells.Item(21, 10).Value = 1.065368697442231
$ws.Cells.Item(21, 11).Value = 1.061221275960037
$ws.Cells.Item(21, 12).Value = 1.067812126793814
$ws.Cells.Item(21, 13).Value = 1.077127306307757

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.058062417993864
$ws.Cells.Item(22, 4).Value = 1.057043291978457
$ws.Cells.Item(22, 5).Value = 1.063539348153077
$ws.Cells.Item(22, 6).Value = 1.072812734631216
$ws.Cells.Item(22, 9).Value = 1.053764617388112
$ws.Cells.Item(22, 10).Value = 1.064730752749946
$ws.Cells.Item(22, 11).Value = 1.060672875549771
$ws.Cells.Item(22, 12).Value = 1.067145082171478
$ws.Cells.Item(22, 13).Value = 1.07638497487314

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.058529766268748
$ws.Cells.Item(23, 4).Value = 1.057402152982424
$ws.Cells.Item(23, 5).Value = 1.063960815035653
$ws.Cells.Item(23, 6).Value = 1.073273544243702
$ws.Cells.Item(23, 9).Value = 1.05393592131768
$ws.Cells.Item(23, 10).Value = 1.065069041685888
$ws.Cells.Item(23, 11).Value = 1.060963707347038
$ws.Cells.Item(23, 12).Value = 1.067498780050007
$ws.Cells.Item(23, 13).Value = 1.076778560360274

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.060368555515566
$ws.Cells.Item(24, 4).Value = 1.058813825163714
$ws.Cells.Item(24, 5).Value = 1.065619531504139
$ws.Cells.Item(24, 6).Value = 1.075087620700129
$ws.Cells.Item(24, 9).Value = 1.05460691113349
$ws.Cells.Item(24, 10).Value = 1.066398594582165
$ws.Cells.Item(24, 11).Value = 1.062106158228138
$ws.Cells.Item(24, 12).Value = 1.068889384130315
$ws.Cells.Item(24, 13).Value = 1.078326717127224

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.062500161785686
$ws.Cells.Item(25, 4).Value = 1.060449706932377
$ws.Cells.Item(25, 5).Value = 1.067543295058516
$ws.Cells.Item(25, 6).Value = 1.077192656323184
$ws.Cells.Item(25, 9).Value = 1.055378518403894
$ws.Cells.Item(25, 10).Value = 1.067936857101576
$ws.Cells.Item(25, 11).Value = 1.063426732263352
$ws.Cells.Item(25, 12).Value = 1.070499295699919
$ws.Cells.Item(25, 13).Value = 1.080120540513131
